$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Prepare the new G (总成绩) column while F (currently the old
#    total-score column) still carries its original formatting, then
#    move the "总成绩" header text from F1 to G1 and retitle F1 as
#    the new "lesson03" header.
# ------------------------------------------------------------------
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$ws.Range("F2:F22").Copy() | Out-Null
$ws.Range("G2:G22").PasteSpecial(-4122) | Out-Null    # xlPasteFormats

$ws.Range("G1").Value = $ws.Range("F1").Value2          # "总成绩" -> G1
$ws.Range("F1").Value = "lesson03"

# ------------------------------------------------------------------
# 2. Give the new lesson03 column (F, rows 2-22) the same formatting
#    as the other score column (lesson02 / E), then fill in scores.
# ------------------------------------------------------------------
$ws.Range("E2:E22").Copy() | Out-Null
$ws.Range("F2:F22").PasteSpecial(-4122) | Out-Null    # xlPasteFormats

$ws.Range("F2").Value = 100
$ws.Range("F3").Value = "-"
$ws.Range("F4").Value = 80
$ws.Range("F5").Value = 80
$ws.Range("F6").Value = 40
$ws.Range("F7").Value = 60
$ws.Range("F8").Value = 80
$ws.Range("F9").Value = 40
$ws.Range("F10").Value = 80
$ws.Range("F11").Value = 60
$ws.Range("F12").Value = 80
$ws.Range("F13").Value = 60
$ws.Range("F14").Value = 80
$ws.Range("F15").Value = 80
$ws.Range("F16").Value = "-"
$ws.Range("F17").Value = 100
$ws.Range("F18").Value = 100
$ws.Range("F19").Value = 100
$ws.Range("F20").Value = 100
$ws.Range("F21").Value = "-"
$ws.Range("F22").Value = 100

# ------------------------------------------------------------------
# 3. 总成绩 (G) = lesson02 + lesson03.
# ------------------------------------------------------------------

# Row 2: first (non-shared) formula of the new total column.
$ws.Range("G2").Formula = '=$E2+$F2'

# Row 3: lesson03 is "-" (text) so E3+F3 errors out; kept as its own
# explicit formula (not part of the shared group below).
$ws.Range("G3").Formula = '=$E3+$F3'

# Rows 4-22: shared formula group.
$ws.Range("G4:G22").Formula = '=$E4+$F4'

# Rows where lesson02 or lesson03 is "-" show a literal "-" instead
# of a formula result (mirrors how the source workbook was edited).
$ws.Range("G7").Value = "-"
$ws.Range("G16").Value = "-"
$ws.Range("G21").Value = "-"

$ws.Calculate()
